$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ordnungssystem")
[void]$ws.Rows.Item(31).Select()
$ws.Rows.Item(31).Delete()
